$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("B2").Value = 10372.65132737054
$ws1.Range("E2").Value = 289260.5393052954
$ws1.Range("G2").Value = 80959.25712661834
$ws1.Range("I2").Value = 161710.6685703679
$ws1.Range("L2").Value = 484922.2142001599
$ws1.Range("M2").Value = 105953.7713982
$ws1.Range("N2").Value = 70003.73489578845
$ws1.Range("O2").Value = 69744.89343456978

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 31203.23858116339
$ws2.Range("E2").Value = 170658.5511254234
$ws2.Range("I2").Value = 209080.6134235085
$ws2.Range("L2").Value = 63518.11613148725
$ws2.Range("M2").Value = 68536.72857011756
$ws2.Range("N2").Value = 19285.19160463996
$ws2.Range("O2").Value = 27033.1386905727

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 27543.1755456332
$ws3.Range("B2").Value = 22113.21643273498
$ws3.Range("E2").Value = 114655.4402706629
$ws3.Range("I2").Value = 153866.0861464091
$ws3.Range("M2").Value = 44638.22942194272
$ws3.Range("N2").Value = 39676.88529639924
$ws3.Range("O2").Value = 31311.04369977792

$ws4 = $wb.Worksheets.Item("2040")
$ws4.Range("N2").Value = 1142.580190039942
$ws4.Range("O2").Value = 0

$ws5 = $wb.Worksheets.Item("2045")
$ws5.Range("A2").Value = 29588.33508286276
$ws5.Range("N2").Value = 4347.543515635315
$ws5.Range("O2").Value = 20429.76977394434
